# WODC/Onderzoeksresultaten.xlsx — "further work on WODC"
#
# Semantic changes applied (the rest of the upstream diff is pure
# churn from the file having been re-saved by a newer Excel build —
# fileVersion/rupBuild, window pixel geometry, theme display name,
# x14ac/mc namespace + knownFonts/dyDescent markers, a default
# slicer-style extLst, and sub-pixel column width / default-row-height
# rounding drift — none of which is reachable or meaningful through
# the Excel object model, so it is intentionally left alone here):
#
#   1. Rename sheet "Bronnen" -> "Deelnemers".
#   2. Move the active sheet/selection from "Definities"!A2 to
#      "Analyse"!E4 (the workbook now opens on the Analyse tab).
#   3. Scroll "Beginselen" so row 19 is the top row, with B25 selected.
#   4. Leave "Deelnemers" (ex-"Bronnen") selection at C8 and
#      "Definities" selection at A2, just no longer the active tab.

$wb = $excel.ActiveWorkbook

# 1. Rename the "Bronnen" sheet to "Deelnemers".
$wsDeelnemers = $wb.Worksheets.Item("Bronnen")
$wsDeelnemers.Name = "Deelnemers"

# 3. Update the view/scroll state on "Beginselen" before switching the
#    active tab away from it, so its own sheet stays frontmost while we
#    set it up.
$wsBeginselen = $wb.Worksheets.Item("Beginselen")
$wsBeginselen.Activate()
$wsBeginselen.Application.ActiveWindow.ScrollRow = 19
$wsBeginselen.Range("B25").Select()

# 2. Make "Analyse" the active sheet with E4 selected (this also clears
#    tabSelected on "Definities", which previously held it).
$wsAnalyse = $wb.Worksheets.Item("Analyse")
$wsAnalyse.Activate()
$wsAnalyse.Range("E4").Select()
